$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Name (row 4) was blank -> now set to "PaysVs" (type correction generated from fsh)
$ws.Range("B4").Value = "PaysVs"

# Date (row 8) gets refreshed to reflect the regeneration timestamp
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
